# Applies the AHDT1_AHD_NEG.docx content edit:
#   "...approximately 2% with the exception of ASXL1 c.1934dup;p.Gly646Trpfs*12
#    (detection limit ~ 5%-10%) and CEBPA (detection limit ~ 10%)."
#   -> "...approximately 4% with the exception of JAK2 c.1849G>T;p.(Val617Phe)
#    (detection limit ~ 1%)."
#
# wdFindContinue = 1, wdReplaceOne = 1 (we only want the single occurrence)

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$oldText = "approximately 2% with the exception of ASXL1 c.1934dup;p.Gly646Trpfs*12 (detection limit ~ 5%-10%) and CEBPA (detection limit ~ 10%)."
$newText = "approximately 4% with the exception of JAK2 c.1849G>T;p.(Val617Phe) (detection limit ~ 1%)."

$result = $d.Content.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

Write-Output "Replace executed: $result"
